$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifting existing rows down) and populate it
# with the 2019 entry, matching the History/Dates page update.
$ws.Rows.Item(2).Insert()
$ws.Cells.Item(2, 1).Value = 2019
$ws.Cells.Item(2, 2).Value = "Andrés Gomez"

# Update the active selection to B3, matching the saved workbook state.
$ws.Range("B3").Select()
